$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SearchData")

for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Value = "yes"
}

$ws.Cells.Item(9, 7).Value = "Quick"
$ws.Cells.Item(15, 8).Value = "uniform-selectProductSort"

$ws.Columns.Item(7).ColumnWidth = 43.6  # engine quantizes to 1/6-character steps; this yields the closest stored width (44.5) to the target 44.5703125

$ws.Range("A7").Select()
